# Apply weekly data refresh to "Hortaliza, Agricola del Norte S.A. de Arica - Pepino dulce"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 44454
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 19000
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = 19500
$ws.Range("P2").Value = 1083
$ws.Range("D3").Value = 44405
$ws.Range("J3").Value = 140
$ws.Range("K3").Value = 17000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 17500
$ws.Range("P3").Value = 972
$ws.Range("D5").Value = 44221
$ws.Range("H5").Value = 'Cultivar XV región'
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = 5500
$ws.Range("N5").Value = '$/caja 10 kilos'
$ws.Range("O5").Value = 'Región de Arica y Parinacota'
$ws.Range("P5").Value = 550
$ws.Range("Q5").Value = 10
$ws.Range("D6").Value = 44554
$ws.Range("H6").Value = 'Cultivar XV región'
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 5500
$ws.Range("N6").Value = '$/caja 10 kilos'
$ws.Range("O6").Value = 'Región de Arica y Parinacota'
$ws.Range("P6").Value = 550
$ws.Range("Q6").Value = 10
$ws.Range("D7").Value = 44391
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 16000
$ws.Range("M7").Value = 15500
$ws.Range("P7").Value = 861
$ws.Range("D8").Value = 44363
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 140
$ws.Range("D9").Value = 44211
$ws.Range("H9").Value = 'Cultivar XV región'
$ws.Range("K9").Value = 4500
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = 4750
$ws.Range("N9").Value = '$/caja 10 kilos'
$ws.Range("O9").Value = 'Región de Arica y Parinacota'
$ws.Range("P9").Value = 475
$ws.Range("Q9").Value = 10
$ws.Range("D10").Value = 44377
$ws.Range("H10").Value = 'Cultivar IV Región'
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 17000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 17600
$ws.Range("N10").Value = '$/bandeja 18 kilos'
$ws.Range("O10").Value = 'Provincia de Limarí'
$ws.Range("P10").Value = 978
$ws.Range("Q10").Value = 18
$ws.Range("D11").Value = 44526
$ws.Range("K11").Value = 5000
$ws.Range("L11").Value = 5500
$ws.Range("M11").Value = 5250
$ws.Range("P11").Value = 525
$ws.Range("D12").Value = 44526
$ws.Range("J12").Value = 100
$ws.Range("L12").Value = 4500
$ws.Range("M12").Value = 4250
$ws.Range("P12").Value = 425
$ws.Range("D13").Value = 44526
$ws.Range("H13").Value = 'Cultivar XV región'
$ws.Range("I13").Value = 'Tercera'
$ws.Range("J13").Value = 120
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 3500
$ws.Range("M13").Value = 3250
$ws.Range("N13").Value = '$/caja 10 kilos'
$ws.Range("O13").Value = 'Región de Arica y Parinacota'
$ws.Range("P13").Value = 325
$ws.Range("Q13").Value = 10
$ws.Range("D14").Value = 44435
$ws.Range("H14").Value = 'Cultivar IV Región'
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 17000
$ws.Range("L14").Value = 18000
$ws.Range("M14").Value = 17500
$ws.Range("N14").Value = '$/bandeja 18 kilos'
$ws.Range("O14").Value = 'Provincia de Limarí'
$ws.Range("P14").Value = 972
$ws.Range("Q14").Value = 18
$ws.Range("I15").Value = 'Tercera'
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14500
$ws.Range("P15").Value = 806
$ws.Range("D16").Value = 44398
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 17000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 17500
$ws.Range("P16").Value = 972
$ws.Range("D17").Value = 44398
$ws.Range("I17").Value = 'Segunda'
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 15500
$ws.Range("P17").Value = 861
$ws.Range("D18").Value = 44533
$ws.Range("K18").Value = 6000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 6500
$ws.Range("P18").Value = 650
$ws.Range("D19").Value = 44533
$ws.Range("J19").Value = 120
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = 4500
$ws.Range("P19").Value = 450
$ws.Range("D20").Value = 44433
$ws.Range("H20").Value = 'Cultivar IV Región'
$ws.Range("I20").Value = 'Segunda'
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 17000
$ws.Range("L20").Value = 18000
$ws.Range("M20").Value = 17500
$ws.Range("N20").Value = '$/bandeja 18 kilos'
$ws.Range("O20").Value = 'Provincia de Limarí'
$ws.Range("P20").Value = 972
$ws.Range("Q20").Value = 18
$ws.Range("D21").Value = 44433
$ws.Range("H21").Value = 'Cultivar IV Región'
$ws.Range("I21").Value = 'Tercera'
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 14500
$ws.Range("N21").Value = '$/bandeja 18 kilos'
$ws.Range("O21").Value = 'Provincia de Limarí'
$ws.Range("P21").Value = 972
$ws.Range("Q21").Value = 18
